$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.046.28"
$ws.Range("E2").Value = '  +2.39%  '
$ws.Range("D3").Value = "'3.232.93"
$ws.Range("E3").Value = '  +6.16%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").Value = "'578.82"
$ws.Range("E5").Value = '  +4.13%  '
$ws.Range("D6").Value = "'150.87"
$ws.Range("E6").Value = '  +6.45%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = "'3.224.91"
$ws.Range("E8").Value = '  +6.43%  '
$ws.Range("E9").Value = '  +4.73%  '
$ws.Range("D10").Value = "'7.06"
$ws.Range("E10").Value = '  +8.82%  '
$ws.Range("E11").Value = '  +5.43%  '
$ws.Range("D12").Value = "'0.487"
$ws.Range("E12").Value = '  +5.16%  '
$ws.Range("D13").Value = "'38.02"
$ws.Range("E13").Value = '  +5.39%  '
$ws.Range("E14").Value = '  +5.46%  '
$ws.Range("D15").Value = "'3.753.98"
$ws.Range("E15").Value = '  +6.53%  '
$ws.Range("D16").Value = "'66.148.98"
$ws.Range("E16").Value = '  +2.38%  '
$ws.Range("D17").Value = "'539.47"
$ws.Range("E17").Value = '  +10.94%  '
$ws.Range("D18").Value = "'3.243.51"
$ws.Range("E18").Value = '  +6.22%  '
$ws.Range("E19").Value = '  +3.06%  '
$ws.Range("D20").Value = "'7.12"
$ws.Range("E20").Value = '  +6.77%  '
$ws.Range("E21").Value = '  +6.31%  '
$ws.Range("D22").Value = "'0.741"
$ws.Range("E22").Value = '  +8.08%  '
$ws.Range("D23").Value = "'7.77"
$ws.Range("E23").Value = '  +9.05%  '
$ws.Range("D24").Value = "'13.49"
$ws.Range("E24").Value = '  +6.72%  '
$ws.Range("D25").Value = "'81.05"
$ws.Range("E25").Value = '  +2.99%  '
$ws.Range("E26").Value = '  +0.21%  '
$ws.Range("E27").Value = '  +18.98%  '
$ws.Range("D28").Value = "'2.95"
$ws.Range("E28").Value = '  +7.87%  '
$ws.Range("E29").Value = '  +7.46%  '
$ws.Range("E30").Value = '  +6.72%  '
$ws.Range("D31").Value = "'2.73"
$ws.Range("E31").Value = '  +4.30%  '
$ws.Range("E32").Value = '  -0.20%  '
$ws.Range("D33").Value = "'1.17"
$ws.Range("E33").Value = '  +5.38%  '
$ws.Range("D34").Value = "'560.99"
$ws.Range("E34").Value = '  +2.54%  '
$ws.Range("E35").Value = '  +6.69%  '
$ws.Range("D36").Value = "'5.60"
$ws.Range("E36").Value = '  +3.56%  '
$ws.Range("D37").Value = "'0.0455"
$ws.Range("E37").Value = '  +8.55%  '
$ws.Range("D38").Value = "'54.83"
$ws.Range("E38").Value = '  +4.41%  '
$ws.Range("E39").Value = '  +7.05%  '
$ws.Range("E40").Value = '  +6.30%  '
$ws.Range("D41").Value = "'3.188.09"
$ws.Range("E41").Value = '  +9.97%  '
$ws.Range("D42").Value = "'2.88"
$ws.Range("E42").Value = '  +3.73%  '
$ws.Range("E43").Value = '  +3.62%  '
$ws.Range("D44").Value = "'0.286"
$ws.Range("E44").Value = '  +17.12%  '
$ws.Range("D45").Value = "'2.33"
$ws.Range("E45").Value = '  +11.53%  '
$ws.Range("D46").Value = "'26.34"
$ws.Range("E46").Value = '  +6.53%  '
$ws.Range("D48").Value = "'0.0₃0552"
$ws.Range("E48").Value = '  +3.12%  '
$ws.Range("D49").Value = "'124.70"
$ws.Range("E49").Value = '  +3.79%  '
$ws.Range("E50").Value = '  +3.26%  '
$ws.Range("E51").Value = '  +7.52%  '
